$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row data: I (ohm value), K (digikey part#), M (digikey URL) ---
$partNumbers = @{
    2  = "36.5KXBK-ND"
    3  = "37.4KXBK-ND"
    4  = "39.2KXBK-ND"
    5  = "33.2KXBK-ND"
    6  = "57.6KXBK-ND"
    7  = "S62KCACT-ND"
    8  = "59.0KXBK-ND"
    9  = "59.0KXBK-ND"
    10 = "59.0KXBK-ND"
    11 = "64.9KXBK-ND"
    12 = "59.0KXBK-ND"
    13 = "56.2KXBK-ND"
}

$urls = @{
    2  = "http://www.digikey.com/product-detail/en/yageo/MFR-25FBF52-36K5/36.5KXBK-ND/13325"
    3  = "http://www.digikey.com/product-detail/en/yageo/MFR-25FBF52-37K4/37.4KXBK-ND/13327"
    4  = "http://www.digikey.com/product-detail/en/yageo/MFR-25FBF52-39K2/39.2KXBK-ND/13331"
    5  = "http://www.digikey.com/product-detail/en/yageo/MFR-25FBF52-33K2/33.2KXBK-ND/13317"
    6  = "http://www.digikey.com/product-detail/en/yageo/MFR-25FBF52-57K6/57.6KXBK-ND/13426"
    7  = "http://www.digikey.com/product-detail/en/stackpole-electronics-inc/RNMF14FTC62K0/S62KCACT-ND/2617523"
    8  = "http://www.digikey.com/product-detail/en/yageo/MFR-25FBF52-59K/59.0KXBK-ND/13428"
    9  = "http://www.digikey.com/product-detail/en/yageo/MFR-25FBF52-59K/59.0KXBK-ND/13428"
    10 = "http://www.digikey.com/product-detail/en/yageo/MFR-25FBF52-59K/59.0KXBK-ND/13428"
    11 = "http://www.digikey.com/product-detail/en/yageo/MFR-25FBF52-64K9/64.9KXBK-ND/13436"
    12 = "http://www.digikey.com/product-detail/en/yageo/MFR-25FBF52-59K/59.0KXBK-ND/13428"
    13 = "http://www.digikey.com/product-detail/en/yageo/MFR-25FBF52-56K2/56.2KXBK-ND/13424"
}

$resistor = @{
    2 = 36500
    3 = 37400
    4 = 39200
    5 = 33200
    6 = 57600
    7 = 62000
    8 = 59000
    9 = 59000
    10 = 59000
    11 = 64900
    12 = 59000
    13 = 56200
}

# --- New column widths / header (Cost column becomes currency, add Digikey URL column) ---
$ws.Columns.Item(12).ColumnWidth = 6.83203125
$ws.Columns.Item(13).ColumnWidth = 92.1640625

# Header L1 ("Cost") becomes bold + currency formatted
$ws.Cells.Item(1, 12).NumberFormat = '_("$"* #,##0.00_);_("$"* \(#,##0.00\);_("$"* "-"??_);_(@_)'

for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 9).Value = $resistor[$r]
    $ws.Cells.Item($r, 10).Formula = "=(F$r*`$B`$2)*I$r"
    $ws.Cells.Item($r, 11).Value = $partNumbers[$r]
    $ws.Cells.Item($r, 12).Value = 0.1
    $ws.Cells.Item($r, 12).NumberFormat = '_("$"* #,##0.00_);_("$"* \(#,##0.00\);_("$"* "-"??_);_(@_)'
    $ws.Cells.Item($r, 13).Value = $urls[$r]
}

# --- sheetView changes ---
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("M20").Select()

# --- workbook window height tweak ---
$excel.ActiveWindow.Height = $excel.ActiveWindow.Height - 80
